# "Updation to testcases" — refresh the "Data" sheet's 12 result rows
# (rows 2-13) so they describe the newest QC run instead of the old one.
#
# Only three columns actually change per row:
#   A = Result ID            (sequential A0948001 .. A0948012)
#   E = Requested Assay/run  (20220311-Cocci-10733Updt, constant for the batch)
#   T = Cartridge ID         (TestCartridge0733, constant for the batch)
# Everything else in the row (Kit Lot, Sample Matrix, dates, etc.) is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$runName     = "20220311-Cocci-10733Updt"
$cartridgeId = "TestCartridge0733"

for ($i = 0; $i -lt 12; $i++) {
    $row = 2 + $i
    $resultId = "A0948" + ("{0:D3}" -f ($i + 1))

    $ws.Cells.Item($row, 1).Value = $resultId    # Column A - Result ID
    $ws.Cells.Item($row, 5).Value = $runName      # Column E - Requested Assay
    $ws.Cells.Item($row, 20).Value = $cartridgeId # Column T - Cartridge ID
}
